$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Rename the "Approach" Heading1 to "Methods".
# ------------------------------------------------------------------
$d.Content.Find.Execute("Approach", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Methods", 2)

# ------------------------------------------------------------------
# 2. Word relocates the hidden "_GoBack" bookmark to the site of the
#    most recent edit whenever the document is saved. Reproduce that
#    by moving the bookmark from the end of the Fig.1 paragraph to
#    immediately after the "Methods" run we just edited.
#
#    A temporary one-character marker is used as the bookmark anchor
#    (COM range objects collapsed to zero length aren't reliable
#    anchors here) and is removed once the bookmark is in place, so
#    the bookmark ends up sitting on its own between the "Methods"
#    run and the paragraph end, exactly where Word would leave it.
# ------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$afterMethods = $d.Content.Duplicate
$afterMethods.Find.Execute("Methods")
$afterMethods.Collapse(0)
$afterMethods.InsertAfter("@")

$marker = $d.Content.Duplicate
$marker.Find.Execute("@")
$d.Bookmarks.Add("_GoBack", $marker)

$marker2 = $d.Content.Duplicate
$marker2.Find.Execute("@")
$marker2.Delete()

# ------------------------------------------------------------------
# 3. The footer's cached PAGE field result moves from "1" to "2".
#    Editing the single character through the Field's Result
#    Characters collection updates the cached field text in place
#    without disturbing the surrounding field-code runs.
# ------------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$pageField = $footer.Range.Fields.Item(1)
$pageField.Result.Characters.Item(1).Text = "2"

$d.Save()
